$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.54
$wsSummary.Range("B4").Value = -0.45
$wsSummary.Range("B6").Value = 207
$wsSummary.Range("B7").Value = 86
$wsSummary.Range("B9").Value = 41.55

# ---- Strategy Status sheet ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C12").Value = 99.15000000000001
$wsStatus.Range("D12").Value = 21
$wsStatus.Range("E12").Value = -0.85
$wsStatus.Range("F12").Value = -0.85
$wsStatus.Range("G12").Value = 28.57

# ---- All Trades sheet: append row 208 (volatility_scorer close) and row 209 (MarketMaking open) ----
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(208, 1).Value = 207
$wsAll.Cells.Item(208, 2).Value = "'2026-02-17"
$wsAll.Cells.Item(208, 3).Value = "10:12:27"
$wsAll.Cells.Item(208, 4).Value = "volatility_scorer"
$wsAll.Cells.Item(208, 5).Value = "NEUTRAL"
$wsAll.Cells.Item(208, 6).Value = 0.05
$wsAll.Cells.Item(208, 7).Value = 0.059481
$wsAll.Cells.Item(208, 8).Value = "CLOSED"
$wsAll.Cells.Item(208, 9).Value = 18.9626
$wsAll.Cells.Item(208, 10).Value = 0.01
$wsAll.Cells.Item(208, 11).Value = 99.15000000000001
$wsAll.Cells.Item(208, 12).Value = 0
$wsAll.Cells.Item(208, 13).Value = 0
$wsAll.Cells.Item(208, 14).Value = 0.85
$wsAll.Cells.Item(208, 15).Value = "Low vol market (score: inf) - ideal for market making"
$wsAll.Cells.Item(208, 16).Value = "early_exit"
$wsAll.Cells.Item(208, 17).Value = 0.18

$wsAll.Cells.Item(209, 1).Value = 208
$wsAll.Cells.Item(209, 2).Value = "'2026-02-17"
$wsAll.Cells.Item(209, 3).Value = "10:12:27"
$wsAll.Cells.Item(209, 4).Value = "MarketMaking"
$wsAll.Cells.Item(209, 5).Value = "DOWN"
$wsAll.Cells.Item(209, 6).Value = 0.05
$wsAll.Cells.Item(209, 7).Value = "'"
$wsAll.Cells.Item(209, 8).Value = "OPEN"
$wsAll.Cells.Item(209, 9).Value = 0
$wsAll.Cells.Item(209, 10).Value = 0
$wsAll.Cells.Item(209, 11).Value = 100.3976537309161
$wsAll.Cells.Item(209, 12).Value = 0
$wsAll.Cells.Item(209, 13).Value = 0
$wsAll.Cells.Item(209, 14).Value = 0.6
$wsAll.Cells.Item(209, 15).Value = "Normal spread capture: 19600 bps"
$wsAll.Cells.Item(209, 16).Value = "'"
$wsAll.Cells.Item(209, 17).Value = 0

# ---- volatility_scorer sheet: append row 22 (trade close) ----
$wsVol = $wb.Worksheets.Item("volatility_scorer")

$wsVol.Cells.Item(22, 1).Value = 207
$wsVol.Cells.Item(22, 2).Value = "'2026-02-17"
$wsVol.Cells.Item(22, 3).Value = "10:12:27"
$wsVol.Cells.Item(22, 4).Value = "volatility_scorer"
$wsVol.Cells.Item(22, 5).Value = "NEUTRAL"
$wsVol.Cells.Item(22, 6).Value = 0.05
$wsVol.Cells.Item(22, 7).Value = 0.059481
$wsVol.Cells.Item(22, 8).Value = "CLOSED"
$wsVol.Cells.Item(22, 9).Value = 18.9626
$wsVol.Cells.Item(22, 10).Value = 0.01
$wsVol.Cells.Item(22, 11).Value = 99.15000000000001
$wsVol.Cells.Item(22, 12).Value = 0
$wsVol.Cells.Item(22, 13).Value = 0
$wsVol.Cells.Item(22, 14).Value = 0.85
$wsVol.Cells.Item(22, 15).Value = "Low vol market (score: inf) - ideal for market making"
$wsVol.Cells.Item(22, 16).Value = "early_exit"
$wsVol.Cells.Item(22, 17).Value = 0.18

# ---- MarketMaking sheet: append row 188 (trade open) ----
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Cells.Item(188, 1).Value = 208
$wsMM.Cells.Item(188, 2).Value = "'2026-02-17"
$wsMM.Cells.Item(188, 3).Value = "10:12:27"
$wsMM.Cells.Item(188, 4).Value = "MarketMaking"
$wsMM.Cells.Item(188, 5).Value = "DOWN"
$wsMM.Cells.Item(188, 6).Value = 0.05
$wsMM.Cells.Item(188, 7).Value = "'"
$wsMM.Cells.Item(188, 8).Value = "OPEN"
$wsMM.Cells.Item(188, 9).Value = 0
$wsMM.Cells.Item(188, 10).Value = 0
$wsMM.Cells.Item(188, 11).Value = 100.3976537309161
$wsMM.Cells.Item(188, 12).Value = 0
$wsMM.Cells.Item(188, 13).Value = 0
$wsMM.Cells.Item(188, 14).Value = 0.6
$wsMM.Cells.Item(188, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(188, 16).Value = "'"
$wsMM.Cells.Item(188, 17).Value = 0
